$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to remain plain text (matches source workbook's
# inlineStr cells) even when the new value looks numeric (e.g. "2.30"),
# then restore the default "Normal" style so no stray number-format is
# left behind on cells that did not have one originally.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '35.539.18'
Set-TextValue $ws.Range("E2") '  -2.78%  '

Set-TextValue $ws.Range("D3") '1.982.62'
Set-TextValue $ws.Range("E3") '  -3.94%  '

Set-TextValue $ws.Range("E4") '  +0.20%  '

Set-TextValue $ws.Range("D5") '241.54'
Set-TextValue $ws.Range("E5") '  -0.45%  '

Set-TextValue $ws.Range("E6") '  -3.74%  '

Set-TextValue $ws.Range("E7") '  +0.07%  '

Set-TextValue $ws.Range("D8") '56.11'
Set-TextValue $ws.Range("E8") '  +6.25%  '

Set-TextValue $ws.Range("D9") '60.13'
Set-TextValue $ws.Range("E9") '  +1.73%  '

Set-TextValue $ws.Range("D10") '0.358'
Set-TextValue $ws.Range("E10") '  -0.59%  '

Set-TextValue $ws.Range("D11") '0.0728'
Set-TextValue $ws.Range("E11") '  -3.13%  '

Set-TextValue $ws.Range("E12") '  -5.02%  '

Set-TextValue $ws.Range("D13") '0.912'
Set-TextValue $ws.Range("E13") '  -0.09%  '

Set-TextValue $ws.Range("D14") '14.18'
Set-TextValue $ws.Range("E14") '  -3.39%  '

Set-TextValue $ws.Range("D15") '2.274.46'
Set-TextValue $ws.Range("E15") '  -3.75%  '

Set-TextValue $ws.Range("D16") '5.22'
Set-TextValue $ws.Range("E16") '  -3.69%  '

Set-TextValue $ws.Range("D17") '2.001.63'
Set-TextValue $ws.Range("E17") '  -2.95%  '

Set-TextValue $ws.Range("D18") '16.96'
Set-TextValue $ws.Range("E18") '  +3.35%  '

Set-TextValue $ws.Range("D19") '35.498.43'
Set-TextValue $ws.Range("E19") '  -2.73%  '

Set-TextValue $ws.Range("D20") '70.18'
Set-TextValue $ws.Range("E20") '  -2.19%  '

Set-TextValue $ws.Range("D21") '0.0₃0834'
Set-TextValue $ws.Range("E21") '  -3.17%  '

Set-TextValue $ws.Range("E22") '  -2.18%  '

Set-TextValue $ws.Range("D23") '5.04'
Set-TextValue $ws.Range("E23") '  -4.32%  '

Set-TextValue $ws.Range("E24") '  +0.06%  '

Set-TextValue $ws.Range("D25") '2.30'
Set-TextValue $ws.Range("E25") '  -3.16%  '

Set-TextValue $ws.Range("D26") '2.29'
Set-TextValue $ws.Range("E26") '  +7.37%  '

Set-TextValue $ws.Range("D27") '163.50'
Set-TextValue $ws.Range("E27") '  -0.50%  '

Set-TextValue $ws.Range("D28") '9.04'
Set-TextValue $ws.Range("E28") '  -4.30%  '

Set-TextValue $ws.Range("D29") '19.41'
Set-TextValue $ws.Range("E29") '  -5.26%  '

Set-TextValue $ws.Range("D30") '0.118'
Set-TextValue $ws.Range("E30") '  -2.80%  '

Set-TextValue $ws.Range("B31") 'Filecoin'
Set-TextValue $ws.Range("C31") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D31") '4.78'
Set-TextValue $ws.Range("E31") '  -5.68%  '

Set-TextValue $ws.Range("B32") 'ImmutableX'
Set-TextValue $ws.Range("C32") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D32") '1.12'
Set-TextValue $ws.Range("E32") '  -2.13%  '

Set-TextValue $ws.Range("E33") '  -1.94%  '

Set-TextValue $ws.Range("D34") '0.0911'
Set-TextValue $ws.Range("E34") '  +11.18%  '

Set-TextValue $ws.Range("D35") '4.23'
Set-TextValue $ws.Range("E35") '  -7.42%  '

Set-TextValue $ws.Range("E36") '  +0.17%  '

Set-TextValue $ws.Range("D37") '2.31'
Set-TextValue $ws.Range("E37") '  +1.66%  '

Set-TextValue $ws.Range("E38") '  -1.79%  '

Set-TextValue $ws.Range("D39") '4.89'
Set-TextValue $ws.Range("E39") '  +1.05%  '

Set-TextValue $ws.Range("D40") '1.18'
Set-TextValue $ws.Range("E40") '  -5.46%  '

Set-TextValue $ws.Range("D41") '2.81'
Set-TextValue $ws.Range("E41") '  -4.17%  '

Set-TextValue $ws.Range("D42") '0.0208'
Set-TextValue $ws.Range("E42") '  -3.27%  '

Set-TextValue $ws.Range("E43") '  -4.80%  '

Set-TextValue $ws.Range("E44") '  -5.18%  '

Set-TextValue $ws.Range("D45") '90.50'
Set-TextValue $ws.Range("E45") '  -4.08%  '

Set-TextValue $ws.Range("D46") '1.370.13'
Set-TextValue $ws.Range("E46") '  -1.39%  '

Set-TextValue $ws.Range("D47") '7.36'
Set-TextValue $ws.Range("E47") '  -0.99%  '

Set-TextValue $ws.Range("D48") '15.38'
Set-TextValue $ws.Range("E48") '  -1.04%  '

Set-TextValue $ws.Range("E49") '  +0.72%  '

Set-TextValue $ws.Range("D50") '2.25'
Set-TextValue $ws.Range("E50") '  -4.59%  '

Set-TextValue $ws.Range("D51") '45.62'
Set-TextValue $ws.Range("E51") '  +0.53%  '
